$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: refine the timestamp precision (same date, tiny floating point update)
$ws.Cells.Item(16, 1).Value = 45816.39142756945

# Row 17: new price entry appended below the existing data
$ws.Cells.Item(17, 1).Value = 45817.39397020341
$ws.Cells.Item(17, 1).NumberFormat = $ws.Cells.Item(16, 1).NumberFormat
$ws.Cells.Item(17, 2).Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Cells.Item(17, 3).Value = "1Kg"
$ws.Cells.Item(17, 4).Value = "15,41€"
